{"js": "// Update the dragspeed readme for the split _X / _Y drag coefficients,\n// add the default filter value, and clarify the 0.5 m/s threshold is\n// per-axis.\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(searchText) +\n      \" but found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Drag term paragraph: DRAGSPEED_COEFF -> DRAGSPEED_COEFF_X and _Y\nawait replaceOnce(\n  \"The drag term (mu/m), defined by DRAGSPEED_COEFF can be found through a short test flight when ground truth velocities are available (e.g. through optitrack).\",\n  \"The drag term (mu/m), defined by DRAGSPEED_COEFF_X and _Y can be found through a short test flight when ground truth velocities are available (e.g. through optitrack).\"\n);\n\n// 2) Low-pass filter strength: append default value\nawait replaceOnce(\n  \"Set the low-pass filter strength (\\u201cGCS\\u2192Settings\\u2192dragspeed\\u2192filter\\u201d or DRAGSPEED_FILTER in the airframe file) to an appropriate value between 0 and 1.\",\n  \"Set the low-pass filter strength (\\u201cGCS\\u2192Settings\\u2192dragspeed\\u2192filter\\u201d or DRAGSPEED_FILTER in the airframe file) to an appropriate value between 0 and 1 (default: 0.8).\"\n);\n\n// 3) Recalibrate DRAGSPEED_COEFF -> DRAGSPEED_COEFF_X and _Y\nawait replaceOnce(\n  \"If the weight of the drone has changed, recalibrate the DRAGSPEED_COEFF.\",\n  \"If the weight of the drone has changed, recalibrate the DRAGSPEED_COEFF_X and _Y.\"\n);\n\n// 4) Fly the drone around: clarify \"along each axis\"\nawait replaceOnce(\n  \"Fly the drone around (manually or using a flight plan). The drag coefficient is updated while the drone flies faster than 0.5 m/s.\",\n  \"Fly the drone around (manually or using a flight plan). The drag coefficient is updated while the drone flies faster than 0.5 m/s along each axis.\"\n);\n\n// 5) Optional Coeff define: DRAGSPEED_COEFF -> DRAGSPEED_COEFF_X and _Y\nawait replaceOnce(\n  \"(Optional) set the new Coeff value (click to update) as a define in the airframe file (DRAGSPEED_COEFF).\",\n  \"(Optional) set the new Coeff value (click to update) as a define in the airframe file (DRAGSPEED_COEFF_X and _Y).\"\n);\n\n// Note: the \"(DRAGSPEED_FILTER set to 0.8)\" / \"(DRAGSPEED_FILTER set to 0.9)\"\n// captions are textually unchanged upstream (the source diff only re-splits\n// their run into two runs with identical combined text), so there is no\n// visible content edit to make there.\n", "ps1": "# Update the dragspeed readme for the split _X / _Y drag coefficients,\n# add the default filter value, and clarify the 0.5 m/s threshold is\n# per-axis.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Exact($searchText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $searchText\"\n    }\n}\n\n# 1) Drag term paragraph: DRAGSPEED_COEFF -> DRAGSPEED_COEFF_X and _Y\nReplace-Exact `\n    \"The drag term (mu/m), defined by DRAGSPEED_COEFF can be found through a short test flight when ground truth velocities are available (e.g. through optitrack).\" `\n    \"The drag term (mu/m), defined by DRAGSPEED_COEFF_X and _Y can be found through a short test flight when ground truth velocities are available (e.g. through optitrack).\"\n\n# 2) Low-pass filter strength: append default value\nReplace-Exact `\n    \"Set the low-pass filter strength (\u201cGCS\u2192Settings\u2192dragspeed\u2192filter\u201d or DRAGSPEED_FILTER in the airframe file) to an appropriate value between 0 and 1.\" `\n    \"Set the low-pass filter strength (\u201cGCS\u2192Settings\u2192dragspeed\u2192filter\u201d or DRAGSPEED_FILTER in the airframe file) to an appropriate value between 0 and 1 (default: 0.8).\"\n\n# 3) Recalibrate DRAGSPEED_COEFF -> DRAGSPEED_COEFF_X and _Y\nReplace-Exact `\n    \"If the weight of the drone has changed, recalibrate the DRAGSPEED_COEFF.\" `\n    \"If the weight of the drone has changed, recalibrate the DRAGSPEED_COEFF_X and _Y.\"\n\n# 4) Fly the drone around: clarify \"along each axis\"\nReplace-Exact `\n    \"Fly the drone around (manually or using a flight plan). The drag coefficient is updated while the drone flies faster than 0.5 m/s.\" `\n    \"Fly the drone around (manually or using a flight plan). The drag coefficient is updated while the drone flies faster than 0.5 m/s along each axis.\"\n\n# 5) Optional Coeff define: DRAGSPEED_COEFF -> DRAGSPEED_COEFF_X and _Y\nReplace-Exact `\n    \"(Optional) set the new Coeff value (click to update) as a define in the airframe file (DRAGSPEED_COEFF).\" `\n    \"(Optional) set the new Coeff value (click to update) as a define in the airframe file (DRAGSPEED_COEFF_X and _Y).\"\n\n# Note: the \"(DRAGSPEED_FILTER set to 0.8)\" / \"(DRAGSPEED_FILTER set to 0.9)\"\n# captions are textually unchanged upstream (the source diff only re-splits\n# their run into two runs with identical combined text), so there is no\n# visible content edit to make there.\n\nWrite-Output \"done\"\n"}
